# Natmi following Dr Hou advice
# Add new cell-cluster combination rows (ECs / FAPs / sCs cross product)
# to the Bmp6-Bmpr1b ligand-receptor pairs sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Bmp6"
$ws.Range("C2").Value = "Bmpr1b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 19.72450833333333
$ws.Range("H2").Value = 59.173525
$ws.Range("I2").Value = 0.5834853563809828
$ws.Range("J2").Value = 0.5834853563809829
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.510190333333334
$ws.Range("N2").Value = 10.530571
$ws.Range("O2").Value = 0.8478537661184122
$ws.Range("P2").Value = 0.8478537661184122
$ws.Range("Q2").Value = 69.23677848141945
$ws.Range("R2").Value = 623.131006332775
$ws.Range("S2").Value = 0.4947102568825602
$ws.Range("T2").Value = 0.4947102568825603

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Bmp6"
$ws.Range("C3").Value = "Bmpr1b"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 19.72450833333333
$ws.Range("H3").Value = 59.173525
$ws.Range("I3").Value = 0.5834853563809828
$ws.Range("J3").Value = 0.5834853563809829
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.629899
$ws.Range("N3").Value = 1.889697
$ws.Range("O3").Value = 0.1521462338815877
$ws.Range("P3").Value = 0.1521462338815877
$ws.Range("Q3").Value = 12.42444807465833
$ws.Range("R3").Value = 111.820032671925
$ws.Range("S3").Value = 0.08877509949842256
$ws.Range("T3").Value = 0.08877509949842259

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Bmp6"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7684289999999999
$ws.Range("H4").Value = 2.305287
$ws.Range("I4").Value = 0.02273146997336134
$ws.Range("J4").Value = 0.02273146997336134
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 3.510190333333334
$ws.Range("N4").Value = 10.530571
$ws.Range("O4").Value = 0.8478537661184122
$ws.Range("P4").Value = 0.8478537661184122
$ws.Range("Q4").Value = 2.697332047653
$ws.Range("R4").Value = 24.275988428877
$ws.Range("S4").Value = 0.01927296242632202
$ws.Range("T4").Value = 0.01927296242632202

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Bmp6"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.7684289999999999
$ws.Range("H5").Value = 2.305287
$ws.Range("I5").Value = 0.02273146997336134
$ws.Range("J5").Value = 0.02273146997336134
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.629899
$ws.Range("N5").Value = 1.889697
$ws.Range("O5").Value = 0.1521462338815877
$ws.Range("P5").Value = 0.1521462338815877
$ws.Range("Q5").Value = 0.4840326586709999
$ws.Range("R5").Value = 4.356293928038999
$ws.Range("S5").Value = 0.003458507547039323
$ws.Range("T5").Value = 0.003458507547039324

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Bmp6"
$ws.Range("C6").Value = "Bmpr1b"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 13.31169566666667
$ws.Range("H6").Value = 39.935087
$ws.Range("I6").Value = 0.3937831736456558
$ws.Range("J6").Value = 0.3937831736456558
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.510190333333334
$ws.Range("N6").Value = 10.530571
$ws.Range("O6").Value = 0.8478537661184122
$ws.Range("P6").Value = 0.8478537661184122
$ws.Range("Q6").Value = 46.72658544940856
$ws.Range("R6").Value = 420.539269044677
$ws.Range("S6").Value = 0.3338705468095299
$ws.Range("T6").Value = 0.33387054680953

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Bmp6"
$ws.Range("C7").Value = "Bmpr1b"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 13.31169566666667
$ws.Range("H7").Value = 39.935087
$ws.Range("I7").Value = 0.3937831736456558
$ws.Range("J7").Value = 0.3937831736456558
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.629899
$ws.Range("N7").Value = 1.889697
$ws.Range("O7").Value = 0.1521462338815877
$ws.Range("P7").Value = 0.1521462338815877
$ws.Range("Q7").Value = 8.385023788737668
$ws.Range("R7").Value = 75.46521409863901
$ws.Range("S7").Value = 0.05991262683612581
$ws.Range("T7").Value = 0.05991262683612583
